{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change: in the \"Requirements\" bullet's detail paragraph\n//   \"Input data format: fastq and fasta file\"\n// append \", output: gui file and genome browser\" right after \"file\" so the\n// paragraph reads:\n//   \"Input data format: fastq and fasta file, output: gui file and genome browser\"\n\nconst body = context.document.body;\n\n// Locate the unique occurrence of \"fasta file\" (end of the target sentence).\nconst results = body.search(\"fasta file\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"fasta file\" in the document body.');\n}\n\nconst target = results.items[0];\ntarget.insertText(\", output: gui file and genome browser\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Change: in the \"Requirements\" bullet's detail paragraph\n#   \"Input data format: fastq and fasta file\"\n# append \", output: gui file and genome browser\" right after \"file\" so the\n# paragraph reads:\n#   \"Input data format: fastq and fasta file, output: gui file and genome browser\"\n\n$d = $word.ActiveDocument\n\n# Find the target paragraph so the Find/Replace below only touches this one\n# spot (there are other \"file\" occurrences elsewhere in the document).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Input data format\")) {\n        $target = $p.Range\n        break\n    }\n}\nif ($target -eq $null) {\n    throw \"Target paragraph ('Input data format...') not found\"\n}\n\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = \" file\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \" file, output: gui file and genome browser\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 1)\n"}
